# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" note text ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$lines = @(
    "Conversión del día 💰",
    "✅ Dólar paralelo: 68",
    "",
    "Binance",
    "✅ 1000 Bs = 5.95 = 24250.0 pesos",
    "✅ 24250.0 pesos = 5.92 = 959.34 Bs",
    "",
    "Promedio competencia",
    "✅ Tasa pesos: 20",
    "✅ Tasa Bs: 20",
    "✅ % Ganancia: 20%"
)
$newText = [string]::Join("`n", $lines)

$wsHoja1.Range("A1").Value = $newText

# --- tasas: update the exchange-rate cells N10, O10, N12, O12 ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 168
$wsTasas.Range("O10").Value = 4074
$wsTasas.Range("N12").Value = 4095
$wsTasas.Range("O12").Value = 162
